$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-preserving column D updates (values look numeric but must stay text,
# matching the original inlineStr/shared-string cell type).

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '30.182.68'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range("E2").Value = '  -0.42%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.861.87'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Range("E3").Value = '  -0.46%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9992'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '241.78'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = '  +2.79%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.9995'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = '  -0.13%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4719'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = '  +0.47%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '42.73'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = '  -0.56%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.2853'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = '  -0.42%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.06473'
$ws.Cells.Item(10, 4).Style = "Normal"

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '20.72'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = '  -4.95%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.07668'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = '  -4.05%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.857.96'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = '  -0.70%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '93.95'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = '  -3.11%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '5.063'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = '  -1.06%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.6803'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = '  -1.28%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '268.93'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = '  -0.15%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '30.178.69'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = '  -0.52%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.34'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = '  -5.67%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.000007542'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = '  -1.47%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.9992'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = '  -0.14%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '2.106.60'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = '  -0.35%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.9991'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '5.170'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range("E24").Value = '  -1.85%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '6.095'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = '  -2.01%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '9.325'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = '  -0.71%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '166.01'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = '  -0.88%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '18.72'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = '  -0.88%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.882'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = '  -3.50%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.376'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = '  +0.67%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.09854'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = '  -0.19%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.505'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = '  +3.27%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.225'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = '  -2.92%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.996'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = '  -1.75%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.04698'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = '  -0.44%  '

$ws.Range("E36").Value = '  -2.44%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.6854'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = '  -2.47%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.707'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = '  -1.22%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01829'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = '  -2.84%  '

$ws.Range("E40").Value = '  -3.64%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '6.367'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = '  +1.67%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '70.06'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = '  -2.97%  '

$ws.Range("E43").Value = '  -0.09%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.8359'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = '  -0.83%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.887'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = '  -3.80%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '101.99'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = '  -1.07%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.4060'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = '  -2.88%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '9.227'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = '  +1.04%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '926.46'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = '  +0.93%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '6.921'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = '  -2.40%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '34.30'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = '  -0.60%  '
